$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Renumber the existing "Test Scenario" index column (A23:A27) ---
# A23 was blank (just bold header style) -> now starts the numbering at 1
$ws.Range("A23").Value = 1
$ws.Range("A24").Value = 2
$ws.Range("A25").Value = 3
$ws.Range("A26").Value = 4
$ws.Range("A27").Value = 5

# --- Add the new test case as row 28 ---
$ws.Range("A28").Value = 6
$ws.Range("B28").Value = "test_CreateUser"
$ws.Range("C28").Value = "This is to test whether users are able to create a new user"
$ws.Range("D28").Value = "Username: JohnnyDoe`r`nPassword: JohnnyDoe@1`r`nPassword confirmation: JohnnyDoe@1"
$ws.Range("E28").Value = "A new user is created"
$ws.Range("F28").Value = "Case failed"

# Match the wrapped-text formatting used by the "Test Values" column elsewhere in the table
$ws.Range("D28").WrapText = $true

# Row 28 holds a multi-line value, so grow its height the same way row 24 (also multi-line) was sized
$ws.Range("D28").RowHeight = 43.75

# --- Update view state to reflect where the author ended up after adding the row ---
$ws.Range("G28").Select()
$excel.ActiveWindow.ScrollRow = 17
$excel.ActiveWindow.ScrollColumn = 4
